# Update the "as_of_utc" timestamp column (AA) on the "Главные" and
# "Линейные" sheets from 2025-11-25 03:03:21 to 2025-11-25 07:09:29.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-11-25 03:03:21"
$newTimestamp = "2025-11-25 07:09:29"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA = 27
        if ($cell.Value2 -eq $oldTimestamp) {
            $cell.Value = $newTimestamp
        }
    }
}
